# This edit inserts a new weekly record at row 59 of the price list
# (pushing every subsequent record down by one row) and fills in the
# new row with the values for the newly reported week.
#
# Net effect on the sheet:
#   - dimension grows from A1:R267 to A1:R268
#   - row 59 becomes a brand new record (date 44608 / volume 120)
#   - every row that used to be at position r (59 <= r <= 267) is now
#     at position r+1, with row 267's old data ending up in new row 268

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 59, shifting rows 59..267 down to 60..268
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new weekly record
$ws.Cells.Item(59, 1).Value = 3
$ws.Cells.Item(59, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44608
$ws.Cells.Item(59, 5).Value = 5
$ws.Cells.Item(59, 6).Value = 100112039
$ws.Cells.Item(59, 7).Value = "Ciboulette"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 120
$ws.Cells.Item(59, 11).Value = 1500
$ws.Cells.Item(59, 12).Value = 1500
$ws.Cells.Item(59, 13).Value = 1500
$ws.Cells.Item(59, 14).Value = "`$/docena de atados"
$ws.Cells.Item(59, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(59, 16).Value = 500
$ws.Cells.Item(59, 17).Value = 3
$ws.Cells.Item(59, 18).Value = "Hortaliza"
